$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "57.064.80"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "3.264.89"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'397.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'108.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  +4.57%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.620"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "'39.44"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "'0.0953"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.44%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").Value = "3.778.41"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "3.272.95"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "'10.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "56.959.47"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("E21").Value = "  +5.96%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").Value = "'292.83"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").Value = "'74.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'28.20"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").Value = "'7.97"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").Value = "'4.37"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "'7.42"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'11.19"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").Value = "'40.08"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.05%  "
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "'51.34"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'3.47"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "'138.27"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'16.72"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").Value = "'22.30"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'2.21"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.35%  "
$ws.Range("D49").Value = "2.150.29"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "'1.99"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.78%  "
$ws.Range("E51").Value = "  -5.84%  "